# Insert a new data row at row 170 (pushing the existing rows 170:271 down
# to 171:272) and populate it with the new record's values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(170).Insert()

$ws.Cells.Item(170, 1).Value2  = 1
$ws.Cells.Item(170, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(170, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(170, 4).Value2  = 44680
$ws.Cells.Item(170, 5).Value2  = 15
$ws.Cells.Item(170, 6).Value2  = 100114013
$ws.Cells.Item(170, 7).Value2  = "Zanahoria"
$ws.Cells.Item(170, 8).Value2  = "Sin especificar"
$ws.Cells.Item(170, 9).Value2  = "Primera"
$ws.Cells.Item(170, 10).Value2 = 100
$ws.Cells.Item(170, 11).Value2 = 24000
$ws.Cells.Item(170, 12).Value2 = 25000
$ws.Cells.Item(170, 13).Value2 = 24500
$ws.Cells.Item(170, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(170, 15).Value2 = "Provincia de Calama"
$ws.Cells.Item(170, 16).Value2 = 980
$ws.Cells.Item(170, 17).Value2 = 25
$ws.Cells.Item(170, 18).Value2 = "Hortaliza"
